# "Generate Report for Handoff"
#
# Updates the localization-status report:
#   - Status moves from "In Translation" to "Ready for handoff"
#     (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
#   - The "Latest Handoff Datetime" / Overview "Latest HO Xliff Generate
#     Date" timestamps are refreshed to reflect the new handoff.
#   - The columns holding the (now longer) status text are widened so the
#     new text fits.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$ws1.Range("G2").Value = "2016-08-27 18:46:50"
$ws3.Range("H2").Value = "2016-08-27 18:46:50"
$ws2.Range("H2").Value = "2016-08-27 18:46:46"

# --- Widen the status columns to fit "Ready for handoff" ---
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
